$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 'This clip shows Akeelah, a young Black girl, competing in the final round of the Scripps National Spelling Bee.'
$ws.Range("I3").Value = 'Tish and Fonny have been together for a few years but are not married. Fonny is in jail for a crime he didn’t commit, and Tish is visiting to share some news'
$ws.Range("I4").Value = 'Frank, a young airline pilot, is taking his father out to lunch at a fancy restaurant.'
$ws.Range("I5").Value = 'Summer and Tom are sitting on a bench, catching up on each other’s lives. They used to be in a romantic relationship, but Summer is now married to another man.'
$ws.Range("I6").Value = 'Troy and Rose have been married for 18 years but are now having an argument because Rose found out that Troy was seeing another woman.'
$ws.Range("I7").Value = 'Forrest visits the grave of his wife, Jenny, delivering a letter that their son, little Forrest, wrote to her.'
$ws.Range("I8").Value = 'In this clip, Will, a young man, attends his first psychotherapy session with Dr. Sean Maguire. They discuss painting and Dr. Maguire’s wife.'
$ws.Range("I9").Value = 'John Coffey, a prisoner on death row, is executed via electric chair. The prison guards watching have come to know John and know that he did not commit the crimes he is being executed for.'
$ws.Range("I20").Value = 'Chris is hoping to get a job as a stockbroker. He shares a cab with Mr. Twistle, a lead manager and partner for the firm he is hoping to work at.'
$ws.Range("I10").Value = 'King George VI addresses Britain and the British Empire over the radio after England enters into World War II with Nazi Germany. King George has been working to overcome a significant stutter with his speech therapist Lionel Logue, whom he calls Logue.'
$ws.Range("I11").Value = 'Christine McPherson wakes up in the hospital after drinking too much during her first weekend at college. She calls her parent’s home and leaves a voicemail.'
$ws.Range("I12").Value = 'Elle, a young female lawyer, takes over the case of Brooke Windham, who is accused of murdering her husband, Hayworth Windham. Elle is questioning Chutney, the daughter of the deceased man.'
$ws.Range("I13").Value = 'This scene shows a conversation at the dinner table. Frank recently came to live with the family of his sister, Sheryl, after attempting suicide. At the table are Sheryl, her brother Frank, her husband Richard, her children Dwayne and Olive, and Richard’s father Edwin.'
$ws.Range("I14").Value = 'This clip shows a fight between Charlie and Nicole, a married couple going through divorce proceedings and a custody battle for their son Henry. Charlie is a theater director, and Nicole used to act in his plays when they lived in New York.'
$ws.Range("I15").Value = 'At the 1980 Winter Olympics, the United States men’s hockey team is getting ready to play their rival, the Soviet Union. Before the game, the head coach, Herb Brooks, talks to the team.'
$ws.Range("I16").Value = 'Juan teaches a young boy, Chiron, how to swim. The two then have a conversation on the beach and Juan gives Chiron some advice.'
$ws.Range("I17").Value = 'This clip shows a hitman getting gas at a gas station in Texas. He talks with the gas station owner and makes him guess the result of his coin toss.'
$ws.Range("I18").Value = 'A young girl named Hal is getting picked up from summer camp by her father. When she arrives to his house, she sees the caretaker and meets her father’s new girlfriend, Meredith.'
$ws.Range("I19").Value = 'Two hitmen, Jules and Vincent, come to an apartment to retrieve a briefcase for their boss, the gangster Marsellus Wallace, from a business partner named Brett.'
$ws.Range("I21").Value = 'Clarice Starling, a young female FBI agent, meets Dr. Hannibal Lector to talk about recent killings by a man known as Buffalo Bill. Dr. Lector is in prison because he is a murderer himself.'
$ws.Range("I22").Value = 'In this clip, a college-aged couple (Marc and Erica) have a conversation in a bar about Marc’s desire to get into a final club at Harvard. Final clubs are social organizations for Harvard students and alumni.'
$ws.Range("I23").Value = 'Robert Graysmith has been investigating the string of killings in and around San Francisco attributed to the Zodiac Killer. He is visiting the house of Bob Vaughn, a silent film organist who used to work with one of Graysmith’s suspects, Rick Marshall.'

$ws.Range("G14").Select() | Out-Null
